# investment-template.xlsx: add gold transaction entity / rework "Gold Transaction" sheet
# - rename sheet "Gold Transaction Statement" -> "Gold Transaction Summary"
# - make the gold sheet the active tab
# - rebuild the gold sheet's layout/columns/headers and add a totals block

$wb = $excel.ActiveWorkbook

$wsFund = $wb.Worksheets.Item(1)
$wsReport = $wb.Worksheets.Item(2)
$wsGold = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. Rename the gold sheet and make it the active tab (this also flips
#    tabSelected between sheet1 and sheet3 + sets bookViews/activeTab).
# ---------------------------------------------------------------------------
$wsGold.Name = "Gold Transaction Summary"
$wsGold.Activate()

# ---------------------------------------------------------------------------
# 2. Wipe the old gold-sheet layout (12 loosely-typed columns, 3 rows) so we
#    can rebuild it with the new 10-column layout.
# ---------------------------------------------------------------------------
$wsGold.Cells.UnMerge() | Out-Null
$wsGold.Cells.Clear() | Out-Null

# Drop the trailing K:L columns (new sheet only spans A:J).
$wsGold.Range("K1:L1").EntireColumn.Delete() | Out-Null

# Uniform column widths A:J (old sheet had 12 bestFit columns).
$wsGold.Range("A1:J1").ColumnWidth = 13.592447916666666

# ---------------------------------------------------------------------------
# 3. Row 1: merged banner "Gold Transactions" across A1:G1, with H1:J1 kept
#    as separate styled (but empty) cells using a vertical-center-only style.
# ---------------------------------------------------------------------------
$bannerRange = $wsGold.Range("A1:J1")
$bannerRange.NumberFormat = "@"
$bannerRange.HorizontalAlignment = -4108   # xlCenter
$bannerRange.VerticalAlignment = -4108     # xlCenter
$bannerRange.Interior.ColorIndex = 0

$wsGold.Cells.Item(1, 1).Value = "Gold Transactions"

# H1:J1 -> new style: numFmt text, vertical-center only (no horizontal).
$tailHeader = $wsGold.Range("H1:J1")
$tailHeader.HorizontalAlignment = 1        # xlGeneral (drops horizontal centering)

$wsGold.Range("A1:G1").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 4. Row 2 / Row 3: column headers + template placeholders.
# ---------------------------------------------------------------------------
$headerRow = $wsGold.Range("A2:G2")
$headerRow.NumberFormat = "@"
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4108

$placeholderRow = $wsGold.Range("A3:G3")
$placeholderRow.NumberFormat = "@"
$placeholderRow.HorizontalAlignment = -4108
$placeholderRow.VerticalAlignment = -4108

$wsGold.Cells.Item(2, 1).Value = "交易类型"
$wsGold.Cells.Item(2, 2).Value = "购买时间"
$wsGold.Cells.Item(2, 3).Value = "赎回时间"
$wsGold.Cells.Item(2, 4).Value = "交易金额"
$wsGold.Cells.Item(2, 5).Value = "交易金额/克"
$wsGold.Cells.Item(2, 6).Value = "克数"
$wsGold.Cells.Item(2, 7).Value = "交易平台"

$wsGold.Cells.Item(3, 1).Value = "{.type}"
$wsGold.Cells.Item(3, 2).Value = "{.purchaseTime}"
$wsGold.Cells.Item(3, 3).Value = "{.redemptionTime}"
$wsGold.Cells.Item(3, 4).Value = "{.amount}"
$wsGold.Cells.Item(3, 5).Value = "{.pircePerGram}"
$wsGold.Cells.Item(3, 7).Value = "{.tradingPlatform}"

# F3 keeps the row-2-style (no fill applied) unlike the rest of row 3.
$wsGold.Range("F3").NumberFormat = "@"
$wsGold.Range("F3").HorizontalAlignment = -4108
$wsGold.Range("F3").VerticalAlignment = -4108
$wsGold.Cells.Item(3, 6).Value = "{.grams}"

# ---------------------------------------------------------------------------
# 5. Row 5 / Row 6: aggregate totals block (row 4 stays blank/untouched).
# ---------------------------------------------------------------------------
$sumHeaderRow = $wsGold.Range("A5:F5")
$sumHeaderRow.NumberFormat = "@"
$sumHeaderRow.HorizontalAlignment = -4108
$sumHeaderRow.VerticalAlignment = -4108

$sumPlaceholderRow = $wsGold.Range("A6:F6")
$sumPlaceholderRow.NumberFormat = "@"
$sumPlaceholderRow.HorizontalAlignment = -4108
$sumPlaceholderRow.VerticalAlignment = -4108

$wsGold.Cells.Item(5, 1).Value = "累计收益"
$wsGold.Cells.Item(5, 2).Value = "累计收益率"
$wsGold.Cells.Item(5, 3).Value = "累计本金"
$wsGold.Cells.Item(5, 4).Value = "累计金额"
$wsGold.Cells.Item(5, 5).Value = "平均金额/克"
$wsGold.Cells.Item(5, 6).Value = "累计克数"

$wsGold.Cells.Item(6, 1).Value = "{sumProfit}"
$wsGold.Cells.Item(6, 2).Value = "{sumYieldRate}"
$wsGold.Cells.Item(6, 3).Value = "{sumPrincipalAmount}"
$wsGold.Cells.Item(6, 4).Value = "{sumAmount}"
$wsGold.Cells.Item(6, 5).Value = "{avgAmountPerGram}"
$wsGold.Cells.Item(6, 6).Value = "{sumGrams}"

# ---------------------------------------------------------------------------
# 6. Rows 7-14: blank, pre-formatted rows (A:J) reserved for future data.
# ---------------------------------------------------------------------------
$blankBlock = $wsGold.Range("A7:J14")
$blankBlock.NumberFormat = "@"
$blankBlock.HorizontalAlignment = -4108
$blankBlock.VerticalAlignment = -4108

Write-Output "Gold sheet rebuilt"
